$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Text fixes (character offsets are computed against the ORIGINAL
#    text of paragraph 1 before any edits; we apply the right-most
#    (later) edit first so the left-most offsets stay valid).
# ---------------------------------------------------------------------

# 1a. "...gymarb finns..." -> "...gymarbete finns..." (insert "ete" right
#     after "gymarb", position 106 in the original text).
$rEte = $d.Range(106, 106)
$rEte.InsertAfter("ete")

# 1b. "udnersökningen" -> "undersökningen" (swap the "dn" typo to "nd").
$rFix = $d.Range(21, 23)
$rFix.Text = "nd"

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the document to right
#    after the newly typed "ete" (position 109 once "ete" is inserted).
#    Remove the old one first.
# ---------------------------------------------------------------------
$oldMark = $d.Bookmarks("_GoBack")
$oldMark.Delete()

$gbRange = $d.Range(109, 109)
$d.Bookmarks.Add("_GoBack", $gbRange)

# ---------------------------------------------------------------------
# 3. Split the runs so the edited words end up as separate <w:r> runs,
#    mirroring how Word keeps each freshly-typed burst in its own run.
#    Each split is done with a scratch bookmark immediately removed --
#    this carves a run boundary without leaving any formatting behind.
#    Work from the right-most boundary to the left-most so earlier
#    offsets are unaffected.
# ---------------------------------------------------------------------

# boundary after "ete" is already implied by the _GoBack bookmark
# (bookmarks themselves force a run split), but make sure the start
# of "ete" (position 106) is also a boundary.
$s1 = $d.Range(106, 106)
$d.Bookmarks.Add("tmpSplit1", $s1)
$d.Bookmarks("tmpSplit1").Delete()

# boundaries around the new "n" and "d" runs (positions 21, 22, 23).
$s2 = $d.Range(23, 23)
$d.Bookmarks.Add("tmpSplit2", $s2)
$d.Bookmarks("tmpSplit2").Delete()

$s3 = $d.Range(22, 22)
$d.Bookmarks.Add("tmpSplit3", $s3)
$d.Bookmarks("tmpSplit3").Delete()

$s4 = $d.Range(21, 21)
$d.Bookmarks.Add("tmpSplit4", $s4)
$d.Bookmarks("tmpSplit4").Delete()
